$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOME")

# Headers for new columns
$ws.Range("C1").Value = "Water Intake in oz"
$ws.Range("D1").Value = "Sleep Hours"
$ws.Range("E1").Value = "Calorie Intake"
$ws.Range("F1").Value = "Challenge Id "

# Sample data rows
$data = @(
    @("James", 10467, 64, 6, 231, 101),
    @("John", 16378, 84, 7, 329, 102),
    @("Robert", 16733, 57, 8, 436, 103),
    @("Michael", 13255, 36, 9, 452, 104),
    @("William", 16722, 47, 6.5, 322, 105),
    @("David", 10278, 98, 5.5, 544, 106),
    @("Richard", 11722, 84, 5, 435, 107),
    @("Joseph", 16832, 71, 7, 235, 108),
    @("Thomas", 18928, 65, 9, 289, 109),
    @("Charles", 19028, 58, 5, 267, 110),
    @("Christopher", 12563, 91, 7, 489, 111),
    @("Daniel", 17383, 101, 8.5, 654, 112),
    @("Matthew", 16738, 89, 4.5, 742, 113),
    @("Anthony", 13893, 61, 6, 341, 114),
    @("Donald", 10273, 82, 8.5, 546, 115),
    @("Mark", 15839, 49, 7.5, 178, 116),
    @("Paul", 14278, 69, 7, 109, 117),
    @("Steven", 10297, 75, 5, 681, 118),
    @("Andrew", 16382, 81, 6, 201, 119),
    @("Kenneth", 9888, 75, 6.5, 301, 120)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row = $row + 1
}

# Autofit the new columns to match bestFit widths recorded in the diff
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null

# Update the selection on HOME sheet and make it the active/selected tab
$ws.Range("I5").Select()

# Make HOME the active sheet (so tabSelected moves from USER to HOME and activeTab=1 in workbook.xml)
$ws.Activate()
